$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos:" block lists two related courses in rows 24 and 25
# (columns B and C hold the same text). The edit reorders them so the
# "LOM3259 ... (Indicação de Conjunto)" entry now comes first (row 24)
# and the "LOM3234 ... (Requisito)" entry comes second (row 25).

$lom3259 = "LOM3259 -  Materiais e Dispositivos Eletrônicos  (Indicação de Conjunto)`n"
$lom3234 = "LOM3234 -  Óptica Física  (Requisito)`n"

$ws.Range("B24").Value = $lom3259
$ws.Range("C24").Value = $lom3259

$ws.Range("B25").Value = $lom3234
$ws.Range("C25").Value = $lom3234
